$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "A86"
$ws.Range("B4").Value = "A89"

$ws.Range("A5").Select() | Out-Null
